$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 5
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("C6").Value = 10
$ws.Range("B7").Value = 12
$ws.Range("C7").Value = 11
$ws.Range("B8").Value = 17
$ws.Range("C8").Value = 22
$ws.Range("B9").Value = 68
$ws.Range("C9").Value = 55
$ws.Range("B10").Value = 69
$ws.Range("C10").Value = 86
$ws.Range("B12").Value = 34
$ws.Range("C12").Value = 38
$ws.Range("B13").Value = 52
$ws.Range("C13").Value = 65
$ws.Range("B14").Value = 35
$ws.Range("C14").Value = 38
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = 10
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 9
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = 6
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 7
$ws.Range("B20").Value = 3
$ws.Range("C20").Value = 7
$ws.Range("B21").Value = 5
$ws.Range("C22").Value = 3
$ws.Range("C24").Value = 2
$ws.Range("B27").Value = 1
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("B32").Value = 2
$ws.Range("C32").Value = 3
$ws.Range("C34").Value = 1
$ws.Range("B36").Value = 2
$ws.Range("C36").Value = 1
$ws.Range("C37").Value = 1
$ws.Range("B38").Value = 3
$ws.Range("B39").Value = 2
$ws.Range("C39").Value = 3
$ws.Range("B40").Value = 3
$ws.Range("C40").Value = 4
$ws.Range("C42").Value = 3
$ws.Range("C43").Value = 6
$ws.Range("B44").Value = 3
$ws.Range("C44").Value = 5
$ws.Range("C45").Value = 3
$ws.Range("C46").Value = 2
$ws.Range("C48").Value = 2
$ws.Range("B49").Value = 1
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = 1
$ws.Range("C51").Value = 1
$ws.Range("B52").Value = 2
$ws.Range("C52").Value = 0
$ws.Range("B53").Value = 2
$ws.Range("C53").Value = 3
